$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KPI Pre")
$ws.Range("D4").Value = 20
$ws.Range("D5").Value = 20
$ws.Range("D6").Value = 20
Write-Host "done"
